$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.467.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3259"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.126"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07062"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.931"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.93%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.602"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.654.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001049"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06597"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.939"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.460.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.474"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.349"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -16.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.840.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.181"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.031"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.713"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -18.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08468"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  -9.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.239"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06054"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02215"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2068"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.01%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.211"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.65%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.200"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.772"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5632"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.946"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06921"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.81%  "
